# Rename Sheet1 and add the three new sheets in the correct order
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

# --- Sheet1: "Sales vs PO" ---
$ws1.Range("A1").Value = "ds"
$ws1.Range("B1").Value = "y"
$ws1.Range("C1").Value = "Order Week"
$ws1.Range("D1").Value = "PO_Requested_Qty"
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

$data = @(
    @(44941,0,44935,0),
    @(44941,0,44935,0),
    @(44976,0,44970,0),
    @(45074,0,45068,0),
    @(45081,0,45075,0),
    @(45088,1,45082,0),
    @(45256,0,45250,0),
    @(45263,0,45257,0),
    @(45270,0,45264,0),
    @(45277,0,45271,0),
    @(45291,0,45285,0),
    @(45298,0,45292,0),
    @(45305,0,45299,0),
    @(45312,0,45306,0),
    @(45319,0,45313,0),
    @(45326,0,45320,0),
    @(45333,0,45327,0),
    @(45340,0,45334,0),
    @(45347,1,45341,0),
    @(45354,0,45348,0),
    @(45361,0,45355,0),
    @(45368,0,45362,0),
    @(45375,0,45369,0),
    @(45382,0,45376,0),
    @(45396,0,45390,0),
    @(45403,0,45397,0),
    @(45410,0,45404,0),
    @(45417,0,45411,0),
    @(45417,0,45411,0),
    @(45424,0,45418,0),
    @(45431,0,45425,0),
    @(45438,0,45432,0),
    @(45445,0,45439,0),
    @(45452,0,45446,0),
    @(45459,0,45453,0),
    @(45466,0,45460,0),
    @(45473,0,45467,0),
    @(45480,0,45474,0),
    @(45487,0,45481,0),
    @(45494,0,45488,0),
    @(45501,0,45495,0),
    @(45508,0,45502,0),
    @(45515,0,45509,0),
    @(45522,1,45516,0),
    @(45529,1,45523,0),
    @(45536,0,45530,0),
    @(45543,0,45537,0),
    @(45550,0,45544,0),
    @(45557,0,45551,0),
    @(45564,0,45558,0),
    @(45571,0,45565,0),
    @(45578,0,45572,0),
    @(45585,0,45579,0),
    @(45592,0,45586,0),
    @(45599,0,45593,0),
    @(45606,0,45600,0),
    @(45613,0,45607,0),
    @(45620,0,45614,0),
    @(45627,0,45621,0),
    @(45634,0,45628,0),
    @(45641,0,45635,0),
    @(45648,0,45642,0),
    @(45655,0,45649,0)
)

$r = 2
foreach ($row in $data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Sheet2: "Weekly Growth" ---
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"
$ws1.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

# --- Sheet3: "Volume Insights" ---
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"
$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# --- Sheet4: "Prediction Info" ---
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Range("A2").Value = 0
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

# Leave the first sheet as the active one
$ws1.Activate()
